$wb = $excel.ActiveWorkbook

# --- Overview sheet: update the "Ready for handoff" row for fe5756af file (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-56-18 02:56:46"

# --- zh-cn sheet: update Status + Latest Handoff Datetime for fe5756af file (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-18 02:56:38"

# --- de-de sheet: update Status + Latest Handoff Datetime for fe5756af file (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-18 02:56:46"
